$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F11").Value = "Documentatiom "
$ws.Range("F12").Value = "Logic "
$ws.Range("F13").Value = "Logic/design "

$ws.Range("H12").Value = 0.41666666666666669
$ws.Range("I12").Value = 0.66666666666666663
$ws.Range("H13").Value = 0.375
$ws.Range("I13").Value = 0.58333333333333337

[void]$ws.Range("F15").Select()

